$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E18 was stored as an inline/text string "2" -> convert to a real number 2
$ws.Cells.Item(18, 5).Value = 2

# Append new row 19 with order data
$ws.Cells.Item(19, 1).Value = "X7W6"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = "Cà phê đen"
$ws.Cells.Item(19, 5).Value = "'1"
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 25000
$ws.Cells.Item(19, 8).Value = 25000
$ws.Cells.Item(19, 9).Value = "2025-05-29 20:10:41"
$ws.Cells.Item(19, 10).Value = "00020101021138550010A000000727012500069704230111440405059060208QRIBFTTA53037045405250005802VN63041B49"
